# Penalty Reward System (unfinished) - remove some now-obsolete weekly/monthly
# data points from both sheets, shifting later rows up.

$wb = $excel.ActiveWorkbook

# "Weekly Quantity" sheet: drop the 5 rows for weeks 45088.99, 45095.99,
# 45102.99, 45123.99 and 45130.99 (old rows 12-16); everything below shifts up.
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Rows("12:16").Delete()

# "Monthly Trend" sheet: drop the 2 rows for months 45107.99 and 45138.99
# (old rows 5-6); everything below shifts up.
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Rows("5:6").Delete()
